# Add the new "2022-Q3" quarterly sheet and update the "总计" (totals) sheet.
#
# Before: 总计, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3
# After : 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3

$wb = $excel.ActiveWorkbook

$totals = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. Insert a brand-new worksheet right before the current "2022-Q2" tab
#    and rename it to "2022-Q3".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2. Populate the new "2022-Q3" sheet with its fund-holdings data.
# ---------------------------------------------------------------------
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "001798"
$newSheet.Cells.Item(2,3).Value = "泰康新回报灵活配置混合A"
$newSheet.Cells.Item(2,4).Value = "1.03"
$newSheet.Cells.Item(2,5).Value = "80.70"
$newSheet.Cells.Item(2,6).Value = "4.32"
$newSheet.Cells.Item(2,7).Value = "0.0445"
$newSheet.Cells.Item(2,8).Value = 7

$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "015201"
$newSheet.Cells.Item(3,3).Value = "创金合信动态平衡混合C"
$newSheet.Cells.Item(3,4).Value = "0.23"
$newSheet.Cells.Item(3,5).Value = "65.33"
$newSheet.Cells.Item(3,6).Value = "2.65"
$newSheet.Cells.Item(3,7).Value = "0.0061"
$newSheet.Cells.Item(3,8).Value = 10

$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "001799"
$newSheet.Cells.Item(4,3).Value = "泰康新回报灵活配置混合C"
$newSheet.Cells.Item(4,4).Value = "0.14"
$newSheet.Cells.Item(4,5).Value = "80.70"
$newSheet.Cells.Item(4,6).Value = "4.32"
$newSheet.Cells.Item(4,7).Value = "0.0060"
$newSheet.Cells.Item(4,8).Value = 7

$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "015200"
$newSheet.Cells.Item(5,3).Value = "创金合信动态平衡混合A"
$newSheet.Cells.Item(5,4).Value = "0.16"
$newSheet.Cells.Item(5,5).Value = "65.33"
$newSheet.Cells.Item(5,6).Value = "2.65"
$newSheet.Cells.Item(5,7).Value = "0.0042"
$newSheet.Cells.Item(5,8).Value = 10

$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).Value = "005281"
$newSheet.Cells.Item(6,3).Value = "中科沃土转型升级灵活配置混合"
$newSheet.Cells.Item(6,4).Value = "0.10"
$newSheet.Cells.Item(6,5).Value = "57.70"
$newSheet.Cells.Item(6,6).Value = "3.01"
$newSheet.Cells.Item(6,7).Value = "0.0030"
$newSheet.Cells.Item(6,8).Value = 7

# Match the header/style formatting used on the sibling quarter sheets
# (bold, centred, bordered header row - style carried by the "基金代码" row).
$q2Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q2Sheet.Range("A2").Copy()
$newSheet.Range("A2:A6").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Update the "总计" summary sheet: insert a new top data row for
#    2022-Q3 and push the existing quarters down by one row.
# ---------------------------------------------------------------------

# Give the new bottom row (row 6) the same "index column" style as the
# rows above it before filling it in.
$totals.Range("A5").Copy()
$totals.Range("A6").PasteSpecial(-4122)
$totals.Cells.Item(6,1).Value = 4

# Shift the B:D (date / count / value) columns down one row, working from
# the bottom up so we never overwrite a row before it's been read.
for ($r = 5; $r -ge 2; $r--) {
    $dest = $r + 1
    $totals.Cells.Item($dest,2).Value = $totals.Cells.Item($r,2).Value2
    $totals.Cells.Item($dest,3).Value = $totals.Cells.Item($r,3).Value2
    $totals.Cells.Item($dest,4).Value = $totals.Cells.Item($r,4).Value2
}

# Write the new 2022-Q3 summary row at the top of the data.
$totals.Cells.Item(2,2).Value = "2022-Q3"
$totals.Cells.Item(2,3).Value = 5
$totals.Cells.Item(2,4).Value = 0.06
